# The deck's theme (ppt/theme/theme1.xml, "Integral" / "Red Violet" colour
# scheme) is switched over to the stock Office theme's colour palette
# ("Office Theme" / "Office" colour scheme). The font scheme and format
# scheme (fills/lines/effects) are already identical between the two
# themes in this deck, so only the twelve theme colours need to change.
#
# PowerPoint exposes the slide master's theme colours through
# Master.ColorScheme.Colors(index).RGB, using the standard
# ppColorSchemeIndex ordering:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

# Target ("Office") colours, in ppColorSchemeIndex order (1-12).
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgbValue = $r + ($g * 256) + ($b * 65536)

    $colorScheme.Colors($i).RGB = $rgbValue
}
